# Update row 9 (Ano 2025) figures in the faturamento_anual sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = 3819835.6
$ws.Range("C9").Value = 602080.24
$ws.Range("D9").Value = 4421915.84
$ws.Range("E9").Value = 13.61582313606403
$ws.Range("F9").Value = 86.38417686393598
$ws.Range("G9").Value = -41.81194691447892
$ws.Range("H9").Value = -31.01904899208732
$ws.Range("I9").Value = 38572
$ws.Range("J9").Value = 1649
$ws.Range("K9").Value = 40221
$ws.Range("L9").Value = 27814
$ws.Range("M9").Value = 158.9816581577623
$ws.Range("N9").Value = 8.540016119698791
